$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 2962.5  # H43: 2845.4546 -> 2962.5
$ws.Cells.Item(43, 9).Value = 2725  # I43: 2755.5557 -> 2725
$ws.Cells.Item(43, 10).Value = 3200  # J43: 3250 -> 3200
$ws.Cells.Item(43, 11).Value = 2725  # K43: 2755.5557 -> 2725
$ws.Cells.Item(43, 12).Value = 3200  # L43: 3250 -> 3200
$ws.Cells.Item(43, 13).Value = -2656  # M43: -2686.5557 -> -2656
$ws.Cells.Item(43, 14).Value = -3338  # N43: -3388 -> -3338
$ws.Cells.Item(98, 8).Value = 98285.25  # H98: 78808.2 -> 98285.25
$ws.Cells.Item(98, 9).Value = 941  # I98: 913.125 -> 941
$ws.Cells.Item(98, 10).Value = 260525.67  # J98: 390388.5 -> 260525.67
$ws.Cells.Item(98, 11).Value = 941  # K98: 913.125 -> 941
$ws.Cells.Item(98, 12).Value = 260525.67  # L98: 390388.5 -> 260525.67
$ws.Cells.Item(98, 13).Value = 557  # M98: 584.875 -> 557
$ws.Cells.Item(98, 14).Value = -263521.67  # N98: -393384.5 -> -263521.67
$ws.Cells.Item(99, 8).Value = 1727.375  # H99: 1662.5714 -> 1727.375
$ws.Cells.Item(99, 9).Value = 1727.375  # I99: 1662.5714 -> 1727.375
$ws.Cells.Item(99, 11).Value = 5182.125  # K99: 4987.7142 -> 5182.125
$ws.Cells.Item(99, 13).Value = -3684.125  # M99: -3489.7142 -> -3684.125
$ws.Cells.Item(106, 8).Value = 36376.13  # H106: 37925.547 -> 36376.13
$ws.Cells.Item(106, 9).Value = 2410.5715  # I106: 2430.8333 -> 2410.5715
$ws.Cells.Item(106, 11).Value = 2410.5715  # K106: 2430.8333 -> 2410.5715
$ws.Cells.Item(106, 13).Value = -1779.5715  # M106: -1799.8333 -> -1779.5715
$ws.Cells.Item(122, 8).Value = 98285.25  # H122: 78808.2 -> 98285.25
$ws.Cells.Item(122, 9).Value = 941  # I122: 913.125 -> 941
$ws.Cells.Item(122, 10).Value = 260525.67  # J122: 390388.5 -> 260525.67
$ws.Cells.Item(122, 11).Value = 2823  # K122: 2739.375 -> 2823
$ws.Cells.Item(122, 12).Value = 781577.01  # L122: 1171165.5 -> 781577.01
$ws.Cells.Item(122, 13).Value = -373  # M122: -289.375 -> -373
$ws.Cells.Item(122, 14).Value = -786477.01  # N122: -1176065.5 -> -786477.01
$ws.Cells.Item(129, 8).Value = 1315.6428  # H129: 1328.5 -> 1315.6428
$ws.Cells.Item(129, 9).Value = 3357.4  # I129: 4074.25 -> 3357.4
$ws.Cells.Item(129, 10).Value = 871.7826  # J129: 870.875 -> 871.7826
$ws.Cells.Item(129, 11).Value = 10072.2  # K129: 12222.75 -> 10072.2
$ws.Cells.Item(129, 12).Value = 2615.3478  # L129: 2612.625 -> 2615.3478
$ws.Cells.Item(129, 13).Value = -5072.200000000001  # M129: -7222.75 -> -5072.200000000001
$ws.Cells.Item(129, 14).Value = -12615.3478  # N129: -12612.625 -> -12615.3478
$ws.Cells.Item(132, 8).Value = 22374.088  # H132: 19101.092 -> 22374.088
$ws.Cells.Item(132, 9).Value = 3031.4055  # I132: 2494.413 -> 3031.4055
$ws.Cells.Item(132, 10).Value = 101894  # J132: 114589.5 -> 101894
$ws.Cells.Item(132, 11).Value = 9094.216499999999  # K132: 7483.239 -> 9094.216499999999
$ws.Cells.Item(132, 12).Value = 305682  # L132: 343768.5 -> 305682
$ws.Cells.Item(132, 13).Value = -6564.216499999999  # M132: -4953.239 -> -6564.216499999999
$ws.Cells.Item(132, 14).Value = -310742  # N132: -348828.5 -> -310742
$ws.Cells.Item(138, 8).Value = 2354.0715  # H138: 2413.5952 -> 2354.0715
$ws.Cells.Item(138, 9).Value = 1415.9412  # I138: 1460.8125 -> 1415.9412
$ws.Cells.Item(138, 10).Value = 2992  # J138: 2999.923 -> 2992
$ws.Cells.Item(138, 11).Value = 4247.8236  # K138: 4382.4375 -> 4247.8236
$ws.Cells.Item(138, 12).Value = 8976  # L138: 8999.769 -> 8976
$ws.Cells.Item(138, 13).Value = 892.1764000000003  # M138: 757.5625 -> 892.1764000000003
$ws.Cells.Item(138, 14).Value = -19256  # N138: -19279.769 -> -19256

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 12283.451  # H32: 13254.851 -> 12283.451
$ws.Cells.Item(32, 9).Value = 11792.707  # I32: 12684.789 -> 11792.707
$ws.Cells.Item(32, 10).Value = 14295.5  # J32: 15661.777 -> 14295.5
$ws.Cells.Item(32, 11).Value = 11792.707  # K32: 12684.789 -> 11792.707
$ws.Cells.Item(32, 12).Value = 14295.5  # L32: 15661.777 -> 14295.5
$ws.Cells.Item(32, 13).Value = -11505.707  # M32: -12397.789 -> -11505.707
$ws.Cells.Item(32, 14).Value = -14869.5  # N32: -16235.777 -> -14869.5
$ws.Cells.Item(46, 8).Value = 8833.333000000001  # H46: 10079.8 -> 8833.333000000001
$ws.Cells.Item(46, 10).Value = 8750  # J46: 10349.75 -> 8750
$ws.Cells.Item(46, 12).Value = 8750  # L46: 10349.75 -> 8750
$ws.Cells.Item(46, 14).Value = -9388  # N46: -10987.75 -> -9388

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 1816.9286  # H134: 1911 -> 1816.9286
$ws.Cells.Item(134, 9).Value = 1394.6522  # I134: 1470.9048 -> 1394.6522
$ws.Cells.Item(134, 11).Value = 4183.9566  # K134: 4412.7144 -> 4183.9566
$ws.Cells.Item(134, 13).Value = -1648.9566  # M134: -1877.7144 -> -1648.9566

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(38, 8).Value = 1804.2222  # H38: 25975 -> 1804.2222
$ws.Cells.Item(38, 9).Value = 1059.5  # I38: 1950 -> 1059.5
$ws.Cells.Item(38, 10).Value = 2400  # J38: 50000 -> 2400
$ws.Cells.Item(38, 11).Value = 1059.5  # K38: 1950 -> 1059.5
$ws.Cells.Item(38, 12).Value = 2400  # L38: 50000 -> 2400
$ws.Cells.Item(38, 13).Value = -682.5  # M38: -1573 -> -682.5
$ws.Cells.Item(38, 14).Value = -3154  # N38: -50754 -> -3154
$ws.Cells.Item(44, 8).Value = 30432.625  # H44: 34017.625 -> 30432.625
$ws.Cells.Item(44, 9).Value = 1400  # I44: 2000 -> 1400
$ws.Cells.Item(44, 10).Value = 34580.145  # J44: 38591.57 -> 34580.145
$ws.Cells.Item(44, 11).Value = 1400  # K44: 2000 -> 1400
$ws.Cells.Item(44, 12).Value = 34580.145  # L44: 38591.57 -> 34580.145
$ws.Cells.Item(44, 13).Value = -958  # M44: -1558 -> -958
$ws.Cells.Item(44, 14).Value = -35464.145  # N44: -39475.57 -> -35464.145
$ws.Cells.Item(45, 8).Value = 10000  # H45: 8955.666999999999 -> 10000
$ws.Cells.Item(45, 9).Value = 0  # I45: 6867 -> 0
$ws.Cells.Item(45, 11).Value = 0  # K45: 6867 -> 0
$ws.Cells.Item(45, 13).ClearContents()  # M45 was -6274
$ws.Cells.Item(46, 8).Value = 1804.2222  # H46: 25975 -> 1804.2222
$ws.Cells.Item(46, 9).Value = 1059.5  # I46: 1950 -> 1059.5
$ws.Cells.Item(46, 10).Value = 2400  # J46: 50000 -> 2400
$ws.Cells.Item(46, 11).Value = 1059.5  # K46: 1950 -> 1059.5
$ws.Cells.Item(46, 12).Value = 2400  # L46: 50000 -> 2400
$ws.Cells.Item(46, 13).Value = -848.5  # M46: -1739 -> -848.5
$ws.Cells.Item(46, 14).Value = -2822  # N46: -50422 -> -2822
$ws.Cells.Item(51, 8).Value = 55589476  # H51: 83365280 -> 55589476
$ws.Cells.Item(51, 10).Value = 38161  # J51: 38331.8 -> 38161
$ws.Cells.Item(51, 12).Value = 38161  # L51: 38331.8 -> 38161
$ws.Cells.Item(51, 14).Value = -39633  # N51: -39803.8 -> -39633
$ws.Cells.Item(61, 8).Value = 55589476  # H61: 83365280 -> 55589476
$ws.Cells.Item(61, 10).Value = 38161  # J61: 38331.8 -> 38161
$ws.Cells.Item(61, 12).Value = 38161  # L61: 38331.8 -> 38161
$ws.Cells.Item(61, 14).Value = -38857  # N61: -39027.8 -> -38857
$ws.Cells.Item(99, 8).Value = 1914.8334  # H99: 1764.4286 -> 1914.8334
$ws.Cells.Item(99, 9).Value = 1601.7142  # I99: 1548 -> 1601.7142
$ws.Cells.Item(99, 10).Value = 2353.2  # J99: 1980.8572 -> 2353.2
$ws.Cells.Item(99, 11).Value = 1601.7142  # K99: 1548 -> 1601.7142
$ws.Cells.Item(99, 12).Value = 2353.2  # L99: 1980.8572 -> 2353.2
$ws.Cells.Item(99, 13).Value = -103.7141999999999  # M99: -50 -> -103.7141999999999
$ws.Cells.Item(99, 14).Value = -5349.2  # N99: -4976.8572 -> -5349.2
$ws.Cells.Item(126, 8).Value = 1914.8334  # H126: 1764.4286 -> 1914.8334
$ws.Cells.Item(126, 9).Value = 1601.7142  # I126: 1548 -> 1601.7142
$ws.Cells.Item(126, 10).Value = 2353.2  # J126: 1980.8572 -> 2353.2
$ws.Cells.Item(126, 11).Value = 4805.142599999999  # K126: 4644 -> 4805.142599999999
$ws.Cells.Item(126, 12).Value = 7059.599999999999  # L126: 5942.571599999999 -> 7059.599999999999
$ws.Cells.Item(126, 13).Value = -2335.142599999999  # M126: -2174 -> -2335.142599999999
$ws.Cells.Item(126, 14).Value = -11999.6  # N126: -10882.5716 -> -11999.6
$ws.Cells.Item(132, 8).Value = 36312.562  # H132: 37203.875 -> 36312.562
$ws.Cells.Item(132, 9).Value = 1751.7931  # I132: 1790.7858 -> 1751.7931
$ws.Cells.Item(132, 11).Value = 5255.379300000001  # K132: 5372.357400000001 -> 5255.379300000001
$ws.Cells.Item(132, 13).Value = -2725.379300000001  # M132: -2842.357400000001 -> -2725.379300000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 554.7041  # H107: 536.2632 -> 554.7041
$ws.Cells.Item(107, 9).Value = 233.85294  # I107: 213.26471 -> 233.85294
$ws.Cells.Item(107, 10).Value = 725.15625  # J107: 716.2951 -> 725.15625
$ws.Cells.Item(107, 11).Value = 701.55882  # K107: 639.79413 -> 701.55882
$ws.Cells.Item(107, 12).Value = 2175.46875  # L107: 2148.8853 -> 2175.46875
$ws.Cells.Item(107, 13).Value = 1218.44118  # M107: 1280.20587 -> 1218.44118
$ws.Cells.Item(107, 14).Value = -6015.46875  # N107: -5988.8853 -> -6015.46875
$ws.Cells.Item(129, 8).Value = 215673.78  # H129: 188770.31 -> 215673.78
$ws.Cells.Item(129, 9).Value = 375952.38  # I129: 300839.9 -> 375952.38
$ws.Cells.Item(129, 10).Value = 1969  # J129: 1987.6666 -> 1969
$ws.Cells.Item(129, 11).Value = 1127857.14  # K129: 902519.7000000001 -> 1127857.14
$ws.Cells.Item(129, 12).Value = 5907  # L129: 5962.9998 -> 5907
$ws.Cells.Item(129, 13).Value = -1122857.14  # M129: -897519.7000000001 -> -1122857.14
$ws.Cells.Item(129, 14).Value = -15907  # N129: -15962.9998 -> -15907
$ws.Cells.Item(137, 8).Value = 10113.637  # H137: 14126.588 -> 10113.637
$ws.Cells.Item(137, 9).Value = 4464.2856  # I137: 3905.6667 -> 4464.2856
$ws.Cells.Item(137, 10).Value = 20000  # J137: 25625.125 -> 20000
$ws.Cells.Item(137, 11).Value = 13392.8568  # K137: 11717.0001 -> 13392.8568
$ws.Cells.Item(137, 12).Value = 60000  # L137: 76875.375 -> 60000
$ws.Cells.Item(137, 13).Value = -8292.856800000001  # M137: -6617.000100000001 -> -8292.856800000001
$ws.Cells.Item(137, 14).Value = -70200  # N137: -87075.375 -> -70200

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2569.359  # H132: 2555.5 -> 2569.359
$ws.Cells.Item(132, 9).Value = 1817.4166  # I132: 1816.6154 -> 1817.4166
$ws.Cells.Item(132, 10).Value = 3772.4666  # J132: 3927.7144 -> 3772.4666
$ws.Cells.Item(132, 11).Value = 5452.2498  # K132: 5449.8462 -> 5452.2498
$ws.Cells.Item(132, 12).Value = 11317.3998  # L132: 11783.1432 -> 11317.3998
$ws.Cells.Item(132, 13).Value = -2922.2498  # M132: -2919.8462 -> -2922.2498
$ws.Cells.Item(132, 14).Value = -16377.3998  # N132: -16843.1432 -> -16377.3998

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 47622480  # H7: 66670570 -> 47622480
$ws.Cells.Item(7, 9).Value = 90911576  # I7: 111113720 -> 90911576
$ws.Cells.Item(7, 10).Value = 4473.8  # J7: 5839.6665 -> 4473.8
$ws.Cells.Item(7, 11).Value = 90911576  # K7: 111113720 -> 90911576
$ws.Cells.Item(7, 12).Value = 4473.8  # L7: 5839.6665 -> 4473.8
$ws.Cells.Item(7, 13).Value = -90911464  # M7: -111113608 -> -90911464
$ws.Cells.Item(7, 14).Value = -4697.8  # N7: -6063.6665 -> -4697.8
$ws.Cells.Item(61, 8).Value = 1686.579  # H61: 1496.4166 -> 1686.579
$ws.Cells.Item(61, 9).Value = 1585  # I61: 1634.5 -> 1585
$ws.Cells.Item(61, 10).Value = 1860.7142  # J61: 1358.3334 -> 1860.7142
$ws.Cells.Item(61, 11).Value = 1585  # K61: 1634.5 -> 1585
$ws.Cells.Item(61, 12).Value = 1860.7142  # L61: 1358.3334 -> 1860.7142
$ws.Cells.Item(61, 13).Value = -1383  # M61: -1432.5 -> -1383
$ws.Cells.Item(61, 14).Value = -2264.7142  # N61: -1762.3334 -> -2264.7142
$ws.Cells.Item(113, 8).Value = 1686.579  # H113: 1496.4166 -> 1686.579
$ws.Cells.Item(113, 9).Value = 1585  # I113: 1634.5 -> 1585
$ws.Cells.Item(113, 10).Value = 1860.7142  # J113: 1358.3334 -> 1860.7142
$ws.Cells.Item(113, 11).Value = 1585  # K113: 1634.5 -> 1585
$ws.Cells.Item(113, 12).Value = 1860.7142  # L113: 1358.3334 -> 1860.7142
$ws.Cells.Item(113, 13).Value = 585  # M113: 535.5 -> 585
$ws.Cells.Item(113, 14).Value = -6200.7142  # N113: -5698.3334 -> -6200.7142
$ws.Cells.Item(126, 8).Value = 47622480  # H126: 66670570 -> 47622480
$ws.Cells.Item(126, 9).Value = 90911576  # I126: 111113720 -> 90911576
$ws.Cells.Item(126, 10).Value = 4473.8  # J126: 5839.6665 -> 4473.8
$ws.Cells.Item(126, 11).Value = 272734728  # K126: 333341160 -> 272734728
$ws.Cells.Item(126, 12).Value = 13421.4  # L126: 17518.9995 -> 13421.4
$ws.Cells.Item(126, 13).Value = -272732258  # M126: -333338690 -> -272732258
$ws.Cells.Item(126, 14).Value = -18361.4  # N126: -22458.9995 -> -18361.4
$ws.Cells.Item(132, 8).Value = 3688.394  # H132: 3656.3823 -> 3688.394
$ws.Cells.Item(132, 9).Value = 2872.2  # I132: 2859.238 -> 2872.2
$ws.Cells.Item(132, 11).Value = 8616.599999999999  # K132: 8577.714 -> 8616.599999999999
$ws.Cells.Item(132, 13).Value = -6086.599999999999  # M132: -6047.714 -> -6086.599999999999
